# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 09:04"

# --- Country name swaps (rows whose rank/order changed) ---
# Argentina / Ucrania swap
$ws.Range("A37").Value = "Ucrania"
$ws.Range("A38").Value = "Argentina"

# Luxemburgo / Hungria swap
$ws.Range("A83").Value = "Hungria"
$ws.Range("A84").Value = "Luxemburgo"

# Groenlandia / Islas Malvinas swap
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"

# Islas Turcas y Caicos / Santa Sede swap
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 7 - India
$ws.Range("B7").Value = 332783
$ws.Range("D7").Value = 169748
$ws.Range("E7").Value = 153512
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 9523

# Row 18 - Pakistan
$ws.Range("B18").Value = 144478
$ws.Range("C18").Value = 5248
$ws.Range("E18").Value = 88028

# Row 37 - now Ucrania
$ws.Range("B37").Value = 31810
$ws.Range("C37").Value = 656
$ws.Range("D37").Value = 14253
$ws.Range("E37").Value = 16656
$ws.Range("G37").Value = 12
$ws.Range("H37").Value = 901

# Row 38 - now Argentina
$ws.Range("B38").Value = 31577
$ws.Range("D38").Value = 9564
$ws.Range("E38").Value = 21180
$ws.Range("H38").Value = 833

# Row 49 - Israel
$ws.Range("B49").Value = 19121
$ws.Range("C49").Value = 66
$ws.Range("D49").Value = 15389
$ws.Range("E49").Value = 3430
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 302

# Row 51 - Barein
$ws.Range("E51").Value = 5366
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 43

# Row 54 - Armenia
$ws.Range("B54").Value = 17064
$ws.Range("C54").Value = 397
$ws.Range("D54").Value = 6276
$ws.Range("E54").Value = 10503
$ws.Range("G54").Value = 16
$ws.Range("H54").Value = 285

# Row 83 - now Hungria
$ws.Range("B83").Value = 4076
$ws.Range("C83").Value = 7
$ws.Range("D83").Value = 2485
$ws.Range("E83").Value = 1028
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 563

# Row 84 - now Luxemburgo
$ws.Range("B84").Value = 4070
$ws.Range("D84").Value = 3929
$ws.Range("E84").Value = 31
$ws.Range("H84").Value = 110

# Row 130 - Georgia
$ws.Range("B130").Value = 879
$ws.Range("C130").Value = 15
$ws.Range("D130").Value = 704
$ws.Range("E130").Value = 161

# Row 197 - Curazao
$ws.Range("D197").Value = 16
$ws.Range("E197").Value = 5

# Row 208 - now Santa Sede
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209 - now Islas Turcas y Caicos
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
